# "flexible working with classes!"
#
# The sheet holds a demand1/net1/pv1/bat1 connection matrix. The column
# headers (row 1) get a "P_from_" prefix, the row headers (column A) get a
# "P_to_" prefix, and every cell that currently holds the flag value 1 is
# replaced with a descriptive flow name "P_<fromLabel>_<toLabel>" (cells
# that are 0 are left as the number 0).
#
# Cells are written in the same order the new labels first appear so the
# regenerated shared-strings table lines up with the target layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header labels: "<label>" -> "P_from_<label>"
$ws.Range("B1").Value = "P_from_demand1"
$ws.Range("C1").Value = "P_from_net1"
$ws.Range("D1").Value = "P_from_pv1"
$ws.Range("E1").Value = "P_from_bat1"

# Column A header labels: "<label>" -> "P_to_<label>"
$ws.Range("A2").Value = "P_to_demand1"
$ws.Range("A3").Value = "P_to_net1"
$ws.Range("A4").Value = "P_to_pv1"
$ws.Range("A5").Value = "P_to_bat1"

# Flagged connections (old value 1) become named flows "P_<from>_<to>".
# Written in the same order the names first appear left-to-right, top-to-
# bottom across the matrix (net1 row, then pv1 row, then bat1 row) so the
# rebuilt shared-strings table lines up with the source layout.
$ws.Range("C2").Value = "P_net1_demand1"
$ws.Range("C5").Value = "P_net1_bat1"
$ws.Range("D2").Value = "P_pv1_demand1"
$ws.Range("D3").Value = "P_pv1_net1"
$ws.Range("D5").Value = "P_pv1_bat1"
$ws.Range("E2").Value = "P_bat1_demand1"
$ws.Range("E3").Value = "P_bat1_net1"

# All remaining cells (pv1 row, and the zero entries above) keep their
# original numeric 0 value untouched.
